# Updates to gold map docs
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Workbook-level view (window position/size) on the bookViews entry
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.Left = 470
$win.Top = 3220
$win.Width = 28530
$win.Height = 17620

# ---------------------------------------------------------------------------
# 2) "Import Strategy" sheet rebuild
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Import Strategy")

# -- Preserve formatting for the reused "Keys/Examples" mini table (rows 1-7)
#    and the "strategy" section banner (row 11) by copying format+values from
#    their old positions (columns E:G) onto the new positions (columns A:C),
#    and the header band onto A11:D11 before the old cells are cleared.
$ws.Range("E1:G7").Copy($ws.Range("A1:C7"))
$ws.Range("E11:G11").Copy($ws.Range("A11:D11"))
$ws.Range("B11:C11").ClearContents()

# -- Preserve the bold "section title" style (used by old A1/E9/A17) for the
#    new "Import strategies" banner at A9, using A17 ("Daily") as the
#    still-untouched format donor.
$ws.Range("A17").Copy($ws.Range("A9"))

# -- Clear the old "Dependencies" table body + "Daily" banner + flowchart
#    content entirely (content & formats) now that the bits worth keeping
#    have been copied to their new homes. Row 8, 10, 12, 15 etc. were/are
#    always blank so a single broad clear is safe.
$ws.Range("E1:G11").Clear()
$ws.Range("E12:G41").Clear()
$ws.Range("A17").Clear()

# New "Keys / Examples" mini table (A1:C7)
$ws.Range("A1").Value = "Type"
$ws.Range("B1").Value = "Keys"
$ws.Range("C1").Value = "Examples"
$ws.Range("A2").Value = "Holding"
$ws.Range("B2").Value = "CPH"
$ws.Range("C2").Value = "12/345/6789"
$ws.Range("A3").Value = "Holder"
$ws.Range("B3").Value = "CPHS"
$ws.Range("C3").Value = "12/345/6789,12/345/7890"
$ws.Range("A4").Value = "Party"
$ws.Range("B4").Value = "PartyId"
$ws.Range("C4").Value = "C1000001"
$ws.Range("A5").Value = "Herd"
$ws.Range("B5").Value = "CPHH"
$ws.Range("C5").Value = "12/345/6789/01"
$ws.Range("B6").Value = "OwnerPartyIds"
$ws.Range("C6").Value = "C1000001"
$ws.Range("B7").Value = "KeeperPartyIds"
$ws.Range("C7").Value = "C1000001,C1000002"

# Section banner
$ws.Range("A9").Value = "Import strategies"

# Strategy headers
$ws.Range("A11").Value = "Bulk import strategy"
$ws.Range("D11").Value = "Daily change strategy"

# Flow body
$ws.Range("A13").Value = "Holdings"
$ws.Range("D13").Value = "Holding changed"
$ws.Range("E13").Value = "For each CPH"

$ws.Range("A14").Value = " - Distinct CPH numbers extracted"
$ws.Range("D14").Value = " - Distinct CPH numbers extracted"
$ws.Range("E14").Value = " - Repeat Bulk Holdings strategy"

$ws.Range("A16").Value = "For each CPH"
$ws.Range("D16").Value = "Herd changed"
$ws.Range("E16").Value = "For each CPH"

$ws.Range("A17").Value = " - Holding by CPH"
$ws.Range("D17").Value = " - Distinct CPH numbers extracted"
$ws.Range("E17").Value = " - Repeat Bulk Holdings strategy"

$ws.Range("A18").Value = " - Herds by CPH(H)"

$ws.Range("A19").Value = " - Parties"
$ws.Range("D19").Value = "Party changed"

$ws.Range("A20").Value = " - Holder by CPHS"
$ws.Range("D20").Value = " - Has a (matching PartyId) holder record changed too?"

$ws.Range("D21").Value = " - If so, discard (Holder will pick up change)"

$ws.Range("A22").Value = "Aggegate Parties and Holders"
$ws.Range("D22").Value = " - If not, continue"

$ws.Range("A23").Value = " - Sam Party"
$ws.Range("B23").Value = "Use Sam Party as base"
$ws.Range("D23").Value = "        - Find herd where owner or keeper id matches (PartyId) and also changed"

$ws.Range("A24").Value = "    - If Sam Holder (matching PartyId) present"
$ws.Range("B24").Value = "Append inferred Holder role. "
$ws.Range("C24").Value = "Check name, communications & address details."
$ws.Range("D24").Value = "           - If present discard (Her will pick up change)"

$ws.Range("A25").Value = " - Sam Holder"
$ws.Range("D25").Value = "           - If not, use PartyId to find existing related CPH numbers"
$ws.Range("E25").Value = "For each CPH"

$ws.Range("A26").Value = "    - If Holder has no associated Party, then use Holder as base"
$ws.Range("E26").Value = " - Repeat Bulk Holdings strategy"

$ws.Range("A27").Value = "    - Else discard"

$ws.Range("D28").Value = "Holder changed"

$ws.Range("A29").Value = "To Silver as is"
$ws.Range("D29").Value = " - Distinct CPH numbers extracted from "
$ws.Range("E29").Value = "For each CPH"

$ws.Range("A30").Value = "To Gold as is"
$ws.Range("E30").Value = " - Repeat Bulk Holdings strategy"

$ws.Range("D31").Value = " - Check for CPH numbers that have been removed"
$ws.Range("E31").Value = "For each CPH orphan"

$ws.Range("E32").Value = " - Remove site party relationships"

# ---------------------------------------------------------------------------
# 3) Column widths / sheet view for "Import Strategy"
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 53.81640625
$ws.Columns.Item(3).ColumnWidth = 53.90625
$ws.Columns.Item(4).ColumnWidth = 67.08984375

$ws.Range("D35").Select()
